# Add new Caltrain capacities (Caltrain PCEP, Caltrain PCBB) to the
# transitVehicleToCapacity sheet, inserted right after the existing
# "Caltrain" row (row 18), pushing all subsequent rows down by 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("transitVehicleToCapacity")

# Update the existing Caltrain row's capacities.
$ws.Range("B18").Value = 1444
$ws.Range("C18").Value = 1228

# Insert two new blank rows at 19 and 20 (formatting is inherited from
# the row above, i.e. style "1", which matches the target layout).
$ws.Range("A19:A20").EntireRow.Insert()

# Row 19: Caltrain PCEP
$ws.Range("A19").Value = "Caltrain PCEP"
$ws.Range("B19").Value = 1502
$ws.Range("C19").Value = 1276.7
$ws.Range("D19").Value = "Caltrain"
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 0
$ws.Range("G19").Value = 0
$ws.Range("H19").Value = 0

# Row 20: Caltrain PCBB
$ws.Range("A20").Value = "Caltrain PCBB"
$ws.Range("B20").Value = 1841
$ws.Range("C20").Value = 1564.85
$ws.Range("D20").Value = "Caltrain"
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 0
$ws.Range("G20").Value = 0
$ws.Range("H20").Value = 0

# Match the final selection recorded in the saved workbook.
$ws.Range("D19").Select()
